$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 680. This shifts the existing rows
# 680..716 down to 681..717 and grows the used range to A1:R717.
$ws.Rows.Item(680).Insert()

# Populate the newly inserted row 680 with the new record.
$ws.Cells.Item(680, 1).Value = 3
$ws.Cells.Item(680, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(680, 3).Value = "Coquimbo"
$ws.Cells.Item(680, 4).Value = 45267
$ws.Cells.Item(680, 5).Value = 5
$ws.Cells.Item(680, 6).Value = 100112031
$ws.Cells.Item(680, 7).Value = "Poroto verde"
$ws.Cells.Item(680, 8).Value = "Magnum"
$ws.Cells.Item(680, 9).Value = "Primera"
$ws.Cells.Item(680, 10).Value = 70
$ws.Cells.Item(680, 11).Value = 32000
$ws.Cells.Item(680, 12).Value = 33000
$ws.Cells.Item(680, 13).Value = 32500
$ws.Cells.Item(680, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(680, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(680, 16).Value = 1300
$ws.Cells.Item(680, 17).Value = 25
$ws.Cells.Item(680, 18).Value = "Hortaliza"
